$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9176648259162903
$ws.Range("B1").Value = 1.243325352668762
$ws.Range("C1").Value = 2.122494220733643
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 2.082007884979248
